# updated main GSC export data
# Adds 4 new daily rows (2025-12-05 .. 2025-12-08) to the "Chart" sheet,
# continuing the existing "No video indexed"/"Video indexed"/"Impressions"
# series, and updates the "Table" sheet's Validation-Failed video count.

$wb = $excel.ActiveWorkbook

$chart = $wb.Worksheets.Item("Chart")
$table = $wb.Worksheets.Item("Table")

# New date labels to append after the last existing row (row 62 = 2025-12-04).
$newDates = @("2025-12-05", "2025-12-06", "2025-12-07", "2025-12-08")

$startRow = 63
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i

    # Write the date as literal text (not an auto-converted date serial):
    # stage it as a formula-computed string in a scratch cell, copy it, and
    # paste-special *values only* into the target cell. Value-only paste
    # carries the string through verbatim instead of re-running Excel's
    # "looks like a date" input parser.
    $chart.Range("Z1").Formula = "=""" + $newDates[$i] + """"
    $chart.Range("Z1").Copy()
    $chart.Range("A" + $row).PasteSpecial(-4163)

    $chart.Range("B" + $row).Value = 23
    $chart.Range("C" + $row).Value = 1
    $chart.Range("D" + $row).Value = 0
}

# Remove the scratch cell used to stage date text.
$chart.Range("Z1").Clear()

# "Videos" / "Failed" validation count drops from 24 to 23.
$table.Range("C2").Value = 23
